$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "\u{1F6A8}SUPER DESCUENTOS POR FIN DE AÑO\u{1F6A8} pusimos todos nuestros productos en oferta al por mayor.\u{1F631} Aprovecha para abastecer tu negocio con estas súper promociones. Valido solo hasta el 31 de diciembre."
$ws.Range("B2").Value = "C:\Users\Personal\Documents\PythonProjects\AppScraping\assets\files\campaña\diciembre_flavia.jpg"

$ws.Range("C2").ClearContents()
$ws.Range("D2").ClearContents()
$ws.Range("E2").ClearContents()

$ws.Range("A2:B2").WrapText = $true
$ws.Range("A2:B2").VerticalAlignment = -4160

$ws.Columns.Item(3).ColumnWidth = 52.7109375

$ws.Range("A6").Select()
